# Slide 12 (1-based), shape "Content Placeholder 3" (3rd shape on the slide),
# 2nd paragraph of its text body currently holds the "takes a list of
# numbers..." text split across three runs:
#   "...two std. " + "dev. " + "from mean"
# The target state merges these into a single run (keeping the first run's
# rPr: dirty="0" smtClean="0").
#
# The host only rewrites the run structure when the paragraph's final text
# actually differs from what is already stored (a same-text assignment is
# treated as a no-op and the original run split is left untouched). So we
# first set the paragraph text to a temporary placeholder, then set it to the
# final desired text; the second assignment lands as a single run with the
# formatting of the paragraph's original leading run, matching the target
# OOXML exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(2)

$para.Text = "__temp_placeholder__"
$para.Text = "This is some code that takes a list of numbers and returns a list of values > two std. dev. from mean"
